$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 29-32: two new "combo" items (015, 016) each split across
# two product lines (red bottle 029 / blue bottle 032), styled with a
# yellow fill (same fill used elsewhere in the sheet) on top of the
# existing thin-border style used by the rest of the table.

$rows = @(
    @{ Row=29; A="015"; B="Combo 2 chai đỏ 3 chai xanh"; C=1; D="029"; E="Giang's Cao xoa thảo dược 50ml (chai đỏ)";              F=2; G="Chai"; H=170000 },
    @{ Row=30; A="015"; B="Combo 2 chai đỏ 3 chai xanh"; C=1; D="032"; E="Giang's Cao xoa thảo dược Giang's 50ml (chai xanh)";    F=3; G="Chai"; H=150000 },
    @{ Row=31; A="016"; B="combo 3 chai đỏ 2 chai xanh"; C=1; D="029"; E="Giang's Cao xoa thảo dược 50ml (chai đỏ)";              F=3; G="Chai"; H=170000 },
    @{ Row=32; A="016"; B="combo 3 chai đỏ 2 chai xanh"; C=1; D="032"; E="Giang's Cao xoa thảo dược Giang's 50ml (chai xanh)";    F=2; G="Chai"; H=150000 }
)

# Seed the two new cell styles once (quote-prefixed text style for the
# "code-like" columns A/D, and a plain style for the rest), both built
# from the existing bordered styles already present in the sheet plus a
# yellow fill, then reused for every subsequent cell.
$quoteStyleSeeded = $false
$plainStyleSeeded = $false

foreach ($r in $rows) {
    $row = $r.Row

    # --- Column A (quote-prefixed text style) ---
    $target = $ws.Range("A$row")
    if (-not $quoteStyleSeeded) {
        $ws.Range("A3").Copy($target)
        $target.Interior.Color = 65535
        $quoteStyleSeeded = $true
    } else {
        $ws.Range("A29").Copy($target)
    }
    $target.Value = "'" + $r.A

    # --- Column B (plain style) ---
    $target = $ws.Range("B$row")
    if (-not $plainStyleSeeded) {
        $ws.Range("B3").Copy($target)
        $target.Interior.Color = 65535
        $plainStyleSeeded = $true
    } else {
        $ws.Range("B29").Copy($target)
    }
    $target.Value = $r.B

    # --- Column C (plain style, numeric) ---
    $target = $ws.Range("C$row")
    $ws.Range("B29").Copy($target)
    $target.Value = $r.C

    # --- Column D (quote-prefixed text style) ---
    $target = $ws.Range("D$row")
    $ws.Range("A29").Copy($target)
    $target.Value = "'" + $r.D

    # --- Column E (plain style) ---
    $target = $ws.Range("E$row")
    $ws.Range("B29").Copy($target)
    $target.Value = $r.E

    # --- Column F (plain style, numeric) ---
    $target = $ws.Range("F$row")
    $ws.Range("B29").Copy($target)
    $target.Value = $r.F

    # --- Column G (plain style) ---
    $target = $ws.Range("G$row")
    $ws.Range("B29").Copy($target)
    $target.Value = $r.G

    # --- Column H (plain style, numeric) ---
    $target = $ws.Range("H$row")
    $ws.Range("B29").Copy($target)
    $target.Value = $r.H
}

$ws.Range("E31").Select()
